# Applies the "HEA-jupyter" test-results column insertion to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column, shifting existing data (rows 2:12, cols B:E) one column
# to the right; row 1's label in B1 stays put.
$ws.Range("B2:B12").Insert(-4161)

# New blank cell at C1 picks up the same style as B1 (left label cell).
$ws.Range("C1").Style = $ws.Range("B1").Style

# Populate the new "HEA-jupyter" column (B) with its results.
$ws.Range("B2").Value = "HEA-jupyter"
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 4.09
$ws.Range("B5").Value = "[fail]"
$ws.Range("B6").Value = 4.09
$ws.Range("B7").Value = 4.21
$ws.Range("B8").Value = 6.2
$ws.Range("B9").Value = 10.71
$ws.Range("B10").Value = 9.4600000000000009
$ws.Range("B11").Value = 27.85
$ws.Range("B12").Value = 7.95

# Match formatting/styles with the corresponding (now-shifted) "HEA-python" column C.
$ws.Range("B2").Style = $ws.Range("C2").Style
$ws.Range("B3").Style = $ws.Range("C3").Style
$ws.Range("B4").Style = $ws.Range("C4").Style
$ws.Range("B5").Style = $ws.Range("C5").Style
$ws.Range("B6").Style = $ws.Range("C6").Style
$ws.Range("B7").Style = $ws.Range("C7").Style
$ws.Range("B8").Style = $ws.Range("C8").Style
$ws.Range("B9").Style = $ws.Range("C9").Style
$ws.Range("B10").Style = $ws.Range("C10").Style
$ws.Range("B11").Style = $ws.Range("C11").Style
$ws.Range("B12").Style = $ws.Range("C12").Style

# Also fix B4/B6 precision the way the author re-entered the numbers (4.09, not
# the long floating point 4.0999999999999996) for the shifted "HEA-python" column.
$ws.Range("C4").Value = 4.09
$ws.Range("C6").Value = 4.09

$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

$ws.Range("D18").Select()
